# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 keeps the same document type / number, but the underlying
# "database" (shared employee roster) was re-ordered, so the Salario
# Basico / Valor Mora entries below were regenerated for the new set
# of outstanding periods.

# Row 16: JOSE DANIEL MENESES ROJAS, periodo 1801
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73139940"
$ws.Range("D16").Value = "JOSE DANIEL MENESES ROJAS"
$ws.Range("E16").Value = "1801"
$ws.Range("F16").Value = 80000
$ws.Range("G16").Value = 2000000

# Row 17: ANGELICA MARIA GALVIS AYAZO, periodo 1808
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45525440"
$ws.Range("D17").Value = "ANGELICA MARIA GALVIS AYAZO"
$ws.Range("E17").Value = "1808"
$ws.Range("F17").Value = 31249
$ws.Range("G17").Value = 781242

# Row 18: ANGELICA MARIA GALVIS AYAZO, periodo 1807
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "45525440"
$ws.Range("D18").Value = "ANGELICA MARIA GALVIS AYAZO"
$ws.Range("E18").Value = "1807"
$ws.Range("F18").Value = 31249
$ws.Range("G18").Value = 781242

# Row 19: ANGELICA MARIA GALVIS AYAZO, periodo 1806
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45525440"
$ws.Range("D19").Value = "ANGELICA MARIA GALVIS AYAZO"
$ws.Range("E19").Value = "1806"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 781242

# Row 20: ANGELICA MARIA GALVIS AYAZO, periodo 1805
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "45525440"
$ws.Range("D20").Value = "ANGELICA MARIA GALVIS AYAZO"
$ws.Range("E20").Value = "1805"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 781242

# Row 21: ANGELICA MARIA GALVIS AYAZO, periodo 1804
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "45525440"
$ws.Range("D21").Value = "ANGELICA MARIA GALVIS AYAZO"
$ws.Range("E21").Value = "1804"
$ws.Range("F21").Value = 31249
$ws.Range("G21").Value = 781242

# Row 22: ANGELICA MARIA GALVIS AYAZO, periodo 1801
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "45525440"
$ws.Range("D22").Value = "ANGELICA MARIA GALVIS AYAZO"
$ws.Range("E22").Value = "1801"
$ws.Range("F22").Value = 29509
$ws.Range("G22").Value = 781242

# Row 23: ANTONIO JOSE STAMBULIE SAER, periodo 1801
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "73091771"
$ws.Range("D23").Value = "ANTONIO JOSE STAMBULIE SAER"
$ws.Range("E23").Value = "1801"
$ws.Range("F23").Value = 80000
$ws.Range("G23").Value = 2000000
